# ----------------------------------------------------------------------------
# Update "BaoCaoTonKhoNguyenLieu" workbook: refresh report timestamp, rename /
# restock raw-material rows, and append new Kiem Ke (stocktake) and Huy Hang
# (disposal) history rows, per commit "cap nhat ttcn va dat ban".
# ----------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # Tong Quan Kho (overview)
$ws2 = $wb.Worksheets.Item(2)   # Chi Tiet Ton Kho (detailed stock)
$ws3 = $wb.Worksheets.Item(3)   # Lich Su Kiem Ke (stocktake diffs)
$ws4 = $wb.Worksheets.Item(4)   # Lich Su Huy Hang (disposal history)

# ---------------------------------------------------------------------------
# Sheet 1 - Tong Quan Kho
# ---------------------------------------------------------------------------

# Report generated-at subtitle.
$ws1.Range("A2").Value = "Báo cáo tạo lúc: 07/11/2025 09:47"

# Key figures.
$ws1.Range("B5").Value = 4168600      # Tong Gia Tri Ton Kho
$ws1.Range("B7").Value = 407500       # Gia Tri Huy (30 ngay qua)

# Column B got a touch narrower.
$ws1.Columns.Item(2).ColumnWidth = 13.75

# ---------------------------------------------------------------------------
# Sheet 2 - Chi Tiet Ton Kho (table grows from 5 to 8 data rows: A3:E8 -> A3:E11)
# ---------------------------------------------------------------------------

# Append 3 new rows (9-11), cloning the look of the first data row (row 4).
foreach ($r in 9..11) {
    $ws2.Range("A4:E4").Copy()
    $ws2.Range("A$r`:E$r").PasteSpecial(-4122)   # xlPasteFormats
}

$ws2.Range("A4").Value  = "Bột Matcha Nhật Bản"
$ws2.Range("B4").Value  = "kg"
$ws2.Range("C4").Value  = 1.94
$ws2.Range("D4").Value  = 0.5
$ws2.Range("E4").Value  = "Đủ dùng"

$ws2.Range("A5").Value  = "Hạt Cà Phê Arabica (Nhập)"
$ws2.Range("B5").Value  = "kg"
$ws2.Range("C5").Value  = 4.76
$ws2.Range("D5").Value  = 1
$ws2.Range("E5").Value  = "Đủ dùng"

$ws2.Range("A6").Value  = "Siro Đào"
$ws2.Range("B6").Value  = "chai"
$ws2.Range("C6").Value  = 8
$ws2.Range("D6").Value  = 2
$ws2.Range("E6").Value  = "Đủ dùng"

$ws2.Range("A7").Value  = "Hạt Cà Phê Robusta (VN)"
$ws2.Range("B7").Value  = "kg"
$ws2.Range("C7").Value  = 10
$ws2.Range("D7").Value  = 2
$ws2.Range("E7").Value  = "Đủ dùng"

$ws2.Range("A8").Value  = "Sữa tươi thanh trùng"
$ws2.Range("B8").Value  = "lít"
$ws2.Range("C8").Value  = 10.2
$ws2.Range("D8").Value  = 5
$ws2.Range("E8").Value  = "Đủ dùng"

$ws2.Range("A9").Value  = "Đường cát trắng"
$ws2.Range("B9").Value  = "kg"
$ws2.Range("C9").Value  = 19.86
$ws2.Range("D9").Value  = 5
$ws2.Range("E9").Value  = "Đủ dùng"

$ws2.Range("A10").Value = "Sữa đặc Ông Thọ (lon)"
$ws2.Range("B10").Value = "lon"
$ws2.Range("C10").Value = 48.79
$ws2.Range("D10").Value = 10
$ws2.Range("E10").Value = "Đủ dùng"

$ws2.Range("A11").Value = "Trà túi lọc Lipton"
$ws2.Range("B11").Value = "túi"
$ws2.Range("C11").Value = 100
$ws2.Range("D11").Value = 20
$ws2.Range("E11").Value = "Đủ dùng"

# Grow the table (Table1) to cover the new rows; autoFilter follows automatically.
$tbl1 = $ws2.ListObjects.Item(1)
$tbl1.Resize($ws2.Range("A3:E11"))

# Extend the stock-level conditional formatting to the new rows.
$fcs2 = $ws2.Range("A4:E8").FormatConditions
for ($i = 1; $i -le $fcs2.Count; $i++) {
    $fcs2.Item($i).ModifyAppliesToRange($ws2.Range("A4:E11"))
}

# Column A widened slightly to fit the longer material names.
$ws2.Columns.Item(1).ColumnWidth = 25.45

# ---------------------------------------------------------------------------
# Sheet 3 - Lich Su Kiem Ke (table grows from 2 to 3 data rows: A3:F5 -> A3:F6)
# ---------------------------------------------------------------------------

# Append 1 new row (6), cloning the look of the first data row (row 4).
$ws3.Range("A4:F4").Copy()
$ws3.Range("A6:F6").PasteSpecial(-4122)   # xlPasteFormats

$ws3.Range("A4").Value = 45964
$ws3.Range("B4").Value = "Trà túi lọc Lipton"
$ws3.Range("C4").Value = 100
$ws3.Range("D4").Value = 98
$ws3.Range("E4").Value = -2
$ws3.Range("F4").Value = "Thất lạc 2 túi"

$ws3.Range("A5").Value = 45961
$ws3.Range("B5").Value = "Sữa tươi thanh trùng"
$ws3.Range("C5").Value = 12
$ws3.Range("D5").Value = 11
$ws3.Range("E5").Value = -1
$ws3.Range("F5").Value = "Hết hạn 1 lít"

$ws3.Range("A6").Value = 45960
$ws3.Range("B6").Value = "Hạt Cà Phê Robusta (VN)"
$ws3.Range("C6").Value = 0.5
$ws3.Range("D6").Value = 0.45
$ws3.Range("E6").Value = -0.05
$ws3.Range("F6").Value = "Hao hụt pha chế"

# Grow the table (Table2) to cover the new row.
$tbl2 = $ws3.ListObjects.Item(1)
$tbl2.Resize($ws3.Range("A3:F6"))

# Extend the over/under conditional formatting to the new row.
$fcs3 = $ws3.Range("E4:E5").FormatConditions
for ($i = 1; $i -le $fcs3.Count; $i++) {
    $fcs3.Item($i).ModifyAppliesToRange($ws3.Range("E4:E6"))
}

# Column B widened slightly to fit the longer material names.
$ws3.Columns.Item(2).ColumnWidth = 23.95

# ---------------------------------------------------------------------------
# Sheet 4 - Lich Su Huy Hang (table grows from 3 to 8 data rows: A3:E6 -> A3:E11)
# ---------------------------------------------------------------------------

# Append 5 new rows (7-11), cloning the look of the first data row (row 4).
foreach ($r in 7..11) {
    $ws4.Range("A4:E4").Copy()
    $ws4.Range("A$r`:E$r").PasteSpecial(-4122)   # xlPasteFormats
}

$ws4.Range("A4").Value  = 45968
$ws4.Range("B4").Value  = "Đường cát trắng"
$ws4.Range("C4").Value  = 2
$ws4.Range("D4").Value  = 20000
$ws4.Range("E4").Value  = "Hủy 2kg đường bị ướt"

$ws4.Range("A5").Value  = 45967
$ws4.Range("B5").Value  = "Sữa đặc Ông Thọ (lon)"
$ws4.Range("C5").Value  = 1
$ws4.Range("D5").Value  = 30000
$ws4.Range("E5").Value  = "Hủy 1 lon sữa đặc móp"

$ws4.Range("A6").Value  = 45966
$ws4.Range("B6").Value  = "Bột Matcha Nhật Bản"
$ws4.Range("C6").Value  = 0.1
$ws4.Range("D6").Value  = 40000
$ws4.Range("E6").Value  = "Hủy 0.1kg bột matcha ẩm mốc"

$ws4.Range("A7").Value  = 45965
$ws4.Range("B7").Value  = "Siro Đào"
$ws4.Range("C7").Value  = 1
$ws4.Range("D7").Value  = 120000
$ws4.Range("E7").Value  = "Pha chế báo hỏng 1 chai siro đào"

$ws4.Range("A8").Value  = 45964
$ws4.Range("B8").Value  = "Hạt Cà Phê Robusta (VN)"
$ws4.Range("C8").Value  = 0.05
$ws4.Range("D8").Value  = 12500
$ws4.Range("E8").Value  = "Hủy 0.05kg cafe hao hụt (sau kiểm kho 1)"

$ws4.Range("A9").Value  = 45963
$ws4.Range("B9").Value  = "Trà túi lọc Lipton"
$ws4.Range("C9").Value  = 2
$ws4.Range("D9").Value  = 10000
$ws4.Range("E9").Value  = "Hủy 2 túi trà (sau kiểm kho 5)"

$ws4.Range("A10").Value = 45962
$ws4.Range("B10").Value = "Hạt Cà Phê Robusta (VN)"
$ws4.Range("C10").Value = 0.5
$ws4.Range("D10").Value = 125000
$ws4.Range("E10").Value = "Hủy nguyên liệu hỏng do trời mưa"

$ws4.Range("A11").Value = 45961
$ws4.Range("B11").Value = "Sữa tươi thanh trùng"
$ws4.Range("C11").Value = 1
$ws4.Range("D11").Value = 50000
$ws4.Range("E11").Value = "Hủy 1 lít sữa tươi hết hạn (sau kiểm kho 2)"

# Grow the table (Table3) to cover the new rows.
$tbl3 = $ws4.ListObjects.Item(1)
$tbl3.Resize($ws4.Range("A3:E11"))

# Extend the "huy > 0" conditional formatting to the new rows.
$fcs4 = $ws4.Range("D4:D6").FormatConditions
for ($i = 1; $i -le $fcs4.Count; $i++) {
    $fcs4.Item($i).ModifyAppliesToRange($ws4.Range("D4:D11"))
}

# Column B widened slightly and column E widened a lot (longer reasons).
$ws4.Columns.Item(2).ColumnWidth = 23.95
$ws4.Columns.Item(5).ColumnWidth = 40.6
